$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank columns before column D (shifts D:K -> F:M)
$ws.Range("D1:E1").EntireColumn.Insert()

# Copy number/date formatting from column F into the two new columns
# for each of the three data blocks, so the new cells inherit the same
# style as the rest of the row (date style for header rows, number style
# for data rows) without generating new style entries.
$ws.Range("F7:F35").Copy()
$ws.Range("D7:E35").PasteSpecial(-4122)
$ws.Range("F38:F77").Copy()
$ws.Range("D38:E77").PasteSpecial(-4122)
$ws.Range("F80:F102").Copy()
$ws.Range("D80:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the two new columns (D = Dec-18, E = Sep-18) with values
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = "NA"
$ws.Range("E8").Value = 97400
$ws.Range("D9").Value = 46100
$ws.Range("E9").Value = 40400
$ws.Range("D10").Value = "NA"
$ws.Range("E10").Value = 57000
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 55000
$ws.Range("E17").Value = 49700
$ws.Range("D18").Value = "NA"
$ws.Range("E18").Value = 47700
$ws.Range("D20").Value = "NA"
$ws.Range("E20").Value = 0
$ws.Range("D21").Value = "NA"
$ws.Range("E21").Value = "NA"
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("D23").Value = -212000
$ws.Range("E23").Value = 47700
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 0
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = -212000
$ws.Range("E26").Value = 47700
$ws.Range("D27").Value = -216200
$ws.Range("E27").Value = 43400
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = "NA"
$ws.Range("E32").Value = 0
$ws.Range("D33").Value = -216200
$ws.Range("E33").Value = 43400
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = -216200
$ws.Range("E35").Value = 43400
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 221700
$ws.Range("E41").Value = 266700
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 22500
$ws.Range("E43").Value = 22400
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("D45").Value = 1900
$ws.Range("E45").Value = 1500
$ws.Range("D46").Value = 0
$ws.Range("E46").Value = 0
$ws.Range("D47").Value = 7991100
$ws.Range("E47").Value = 8029400
$ws.Range("D48").Value = 0
$ws.Range("E48").Value = 0
$ws.Range("D49").Value = 0
$ws.Range("E49").Value = 0
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 10500
$ws.Range("E52").Value = 15100
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 8464600
$ws.Range("E54").Value = 8680900
$ws.Range("D57").Value = 3600
$ws.Range("E57").Value = 6300
$ws.Range("D58").Value = 7037700
$ws.Range("E58").Value = 7218700
$ws.Range("D59").Value = 107500
$ws.Range("E59").Value = 241700
$ws.Range("D60").Value = 0
$ws.Range("E60").Value = 0
$ws.Range("D61").Value = 0
$ws.Range("E61").Value = 0
$ws.Range("D62").Value = 0
$ws.Range("E62").Value = 0
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 7339300
$ws.Range("E66").Value = 7472600
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = -1583200
$ws.Range("E72").Value = -1342500
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 1125300
$ws.Range("E76").Value = 1208300
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = -216200
$ws.Range("E81").Value = 43400
$ws.Range("D83").Value = 0
$ws.Range("E83").Value = 0
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 27800
$ws.Range("E89").Value = 30100
$ws.Range("D91").Value = 0
$ws.Range("E91").Value = 0
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = 101400
$ws.Range("E94").Value = -1044000
$ws.Range("D96").Value = -24500
$ws.Range("E96").Value = -24200
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -180900
$ws.Range("E100").Value = 945100
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = -51700
$ws.Range("E102").Value = -68900
